# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets,
# matching the data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 6615
$sheet1.Range("F4").Value = 413
$sheet1.Range("F7").Value = 540
$sheet1.Range("F15").Value = 3293
$sheet1.Range("F16").Value = 16
$sheet1.Range("F18").Value = 1927
$sheet1.Range("F19").Value = 48
$sheet1.Range("F21").Value = 124

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 6615
$sheet4.Range("F4").Value = 413
$sheet4.Range("F8").Value = 540
$sheet4.Range("F16").Value = 3293
$sheet4.Range("F17").Value = 16
$sheet4.Range("F19").Value = 1927
$sheet4.Range("F20").Value = 48
$sheet4.Range("F22").Value = 124
